$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "Player Info" sheet before the existing "ODI Batting"
#    sheet (Worksheets.Add() inserts before the active sheet, which is
#    "ODI Batting" on load, so this lands it first).
# ------------------------------------------------------------------
$info = $wb.Worksheets.Add()
$info.Name = "Player Info"

# Re-fetch the "ODI Batting" handle *after* the insert -- sheet
# references are positional, so a handle grabbed before the insert
# would now (incorrectly) point at the new sheet.
$odi = $wb.Worksheets.Item("ODI Batting")

# Header row (bold, bordered, centered like the existing sheet headers)
$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $headers.Length; $c++) {
    $cell = $info.Cells.Item(1, $c)
    $cell.Value = $headers[$c - 1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Data row — ID stored as text (matches the source data's inline-string type)
$info.Cells.Item(2, 1).NumberFormat = "@"
$info.Cells.Item(2, 1).Value = "4651"
$info.Cells.Item(2, 2).Value = "Ihsanullah Janat"
$info.Cells.Item(2, 3).Value = "Right Handed"
$info.Cells.Item(2, 4).Value = "Does Not Bowl | Unknown"

# ------------------------------------------------------------------
# 2. Rework the "ODI Batting" sheet's MATCH_CARD_LINK column into a
#    MATCH_CODE column -- keep just the numeric code that used to be
#    the `MatchCode=` query-string parameter on the howstat.com URL.
# ------------------------------------------------------------------
$odi.Range("D1").Value = "MATCH_CODE"

$lastRow = $odi.UsedRange.Rows.Count
$odi.Range("D2:D$lastRow").NumberFormat = "@"
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $odi.Cells.Item($r, 4)
    $url = $cell.Value()
    $idx = $url.IndexOf("MatchCode=")
    $code = $url.Substring($idx + 10)
    $cell.Value = $code
}
